# Bugfixed the naive forecaster component module
# Updates the forecast-error table values in B2:G11 to reflect the fixed
# naive forecaster computations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  -0.1052853234940633, 0.4900476106607171, 0.6183260505831405, 0.7863371100127098, 0.7868594147108519, 52),
    @(3,   0.3646584415625562, 0.722637768445621,  1.426035144228675,  1.194167134126825,  1.148442664015213,  51),
    @(4,   0.4532797608695453, 0.9902620268281389, 3.506287014911406,  1.872508214911595,  1.835262443918751,  50),
    @(5,   0.4187269926694491, 1.116030423578108,  4.462820609660299,  2.112538901336564,  2.092082930599334,  49),
    @(6,   0.2603619109577968, 0.935772708153241,  3.863674849690833,  1.965623272575606,  1.96892104161102,   48),
    @(7,   0.2925042579902095, 0.9789877322222572, 4.769865589305868,  2.184002195352804,  2.19502714946493,   36),
    @(8,   0.2341597008279832, 0.9857524095050276, 4.861337304003894,  2.204844054350306,  2.224381796787603,  35),
    @(9,   0.1803491558900733, 1.423780144984279,  8.982173702935924,  2.997027477841324,  3.078327104205846,  18),
    @(10, -0.6058116565806465, 1.074425645603604,  6.535069776881333,  2.556378253874284,  2.60477795886628,   11),
    @(11,  0.1218082912510491, 0.4641168854429287, 0.2644604851195395, 0.5142572168861994, 0.5585955886216332,  5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}
